# Add data for 2022-11-30
# This applies incremental updates (from newly-arrived incident records)
# to the YTD crime-count workbook: a handful of per-neighborhood sheets,
# the "By Neighborhood" roll-up sheet, and the "Citywide Totals" roll-up
# sheet each get a small number of cells bumped to their new totals.

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        [string]$SheetName,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cell in $Values.Keys) {
        $ws.Range($cell).Value = $Values[$cell]
    }
}

Set-Cells "Citywide Totals" @{
    "F2" = 91
    "D3" = 132
    "B6" = 368
    "C6" = 471
    "D6" = 407
    "E6" = 460
    "F6" = 515
    "G6" = 431
    "I6" = 493
    "B7" = 493
    "C7" = 625
    "D7" = 636
    "E7" = 683
    "F7" = 747
    "G7" = 660
    "I7" = 823
}

Set-Cells "By Neighborhood" @{
    "D10" = 3
    "F19" = 23
    "D28" = 45
    "C32" = 39
    "I32" = 48
    "G36" = 27
    "D47" = 14
    "E53" = 81
    "F53" = 79
    "C63" = 7
    "B80" = 15
    "D82" = 10
    "C85" = 15
    "B98" = 493
    "C98" = 625
    "D98" = 636
    "E98" = 683
    "F98" = 747
    "G98" = 660
    "I98" = 823
}

Set-Cells "Garfield Park" @{
    "C6" = 34
    "I6" = 32
    "C7" = 39
    "I7" = 48
}

Set-Cells "Grand Crossing" @{
    "G6" = 15
    "G7" = 27
}

Set-Cells "South Chicago" @{
    "B4" = 12
    "B5" = 15
}

Set-Cells "Englewood" @{
    "D3" = 17
    "D7" = 45
}

Set-Cells "South Shore" @{
    "D4" = 6
    "D5" = 10
}

Set-Cells "Loop" @{
    "F2" = 7
    "E6" = 63
    "F6" = 58
    "E7" = 81
    "F7" = 79
}

Set-Cells "Avondale" @{
    "D5" = 3
    "D6" = 3
}

Set-Cells "New City" @{
    "C4" = 6
    "C5" = 7
}

Set-Cells "Chatham" @{
    "F5" = 16
    "F6" = 23
}

Set-Cells "United Center" @{
    "C4" = 12
    "C5" = 15
}

Set-Cells "Lake View" @{
    "D5" = 12
    "D6" = 14
}
